$d = $word.ActiveDocument

$replacements = @(
    @('0+39=39', '76-43=33'),
    @('29+25=54', '8+47=55'),
    @('74-48=26', '53-49=4'),
    @('52-17=35', '81-18=63'),
    @('95-34=61', '2+35=37'),
    @('53-20=33', '56+19=75'),
    @('46-29=17', '20+73=93'),
    @('57-8=49', '52-46=6'),
    @('65-55=10', '41+18=59'),
    @('43+46=89', '5+14=19'),
    @('42-19=23', '98-23=75'),
    @('73-61=12', '94+2=96'),
    @('63-30=33', '16+46=62'),
    @('57-42=15', '19+4=23'),
    @('16+61=77', '8-2=6'),
    @('25-23=2', '75-43=32'),
    @('37+8=45', '30-0=30'),
    @('83-15=68', '81-16=65'),
    @('5+84=89', '91-9=82'),
    @('7+24=31', '29-19=10'),
    @('26-16=10', '87-47=40'),
    @('96-87=9', '73+4=77'),
    @('81-3=78', '99-96=3'),
    @('17+51=68', '60+26=86'),
    @('59-58=1', '54+40=94'),
    @('63-29=34', '19+9=28'),
    @('29+13=42', '23+25=48'),
    @('59+21=80', '54+32=86'),
    @('46+47=93', '98-50=48'),
    @('15+78=93', '68+2=70'),
    @('67+0=67', '79-7=72'),
    @('78+7=85', '82+11=93'),
    @('88+8=96', '49+41=90'),
    @('40+16=56', '28+30=58'),
    @('32-3=29', '57-18=39'),
    @('82-5=77', '24-18=6'),
    @('63-44=19', '28+17=45'),
    @('49-16=33', '81-0=81'),
    @('54+18=72', '36-11=25'),
    @('46+7=53', '77-59=18'),
    @('44-32=12', '19-6=13'),
    @('22+67=89', '70-37=33'),
    @('42-7=35', '93+4=97'),
    @('22+63=85', '71+18=89'),
    @('26-6=20', '46+11=57'),
    @('89-76=13', '66-18=48'),
    @('10+84=94', '89+5=94'),
    @('23-23=0', '48-13=35'),
    @('69-48=21', '88+10=98'),
    @('18+30=48', '81-25=56'),
    @('92-73=19', '89-63=26'),
    @('13+2=15', '15+54=69'),
    @('72+9=81', '25+47=72'),
    @('56+21=77', '55+44=99'),
    @('95-75=20', '91-88=3'),
    @('49-45=4', '79-40=39'),
    @('25+74=99', '68-66=2'),
    @('13+29=42', '73-17=56'),
    @('64-12=52', '24+32=56'),
    @('17+14=31', '25-11=14'),
    @('19+19=38', '5+62=67'),
    @('96+1=97', '52+18=70'),
    @('79-32=47', '44+14=58'),
    @('22+13=35', '60+16=76'),
    @('45+46=91', '3+46=49'),
    @('94-86=8', '24+22=46'),
    @('31+31=62', '75-36=39'),
    @('11+35=46', '89-80=9'),
    @('43+55=98', '90-2=88'),
    @('14+32=46', '30-0=30'),
    @('3+49=52', '31+57=88'),
    @('72+8=80', '39+58=97'),
    @('27+22=49', '18+12=30'),
    @('26+43=69', '22+30=52'),
    @('97-60=37', '3+19=22'),
    @('10+47=57', '21+76=97'),
    @('54-36=18', '72-25=47'),
    @('22+11=33', '29+28=57'),
    @('62-34=28', '25-22=3'),
    @('24+38=62', '10+64=74'),
    @('8+70=78', '11+65=76'),
    @('5+20=25', '15+42=57'),
    @('88-56=32', '42-20=22'),
    @('92-83=9', '36+53=89'),
    @('30-27=3', '97-52=45'),
    @('14+53=67', '25-15=10'),
    @('13+49=62', '43-3=40'),
    @('57+29=86', '25+32=57'),
    @('58+26=84', '56+18=74'),
    @('74-49=25', '48+17=65'),
    @('63+36=99', '66-53=13'),
    @('89-68=21', '5+82=87'),
    @('89-29=60', '50-26=24'),
    @('94-15=79', '90-88=2'),
    @('78-69=9', '73+3=76'),
    @('20+49=69', '11+67=78'),
    @('47+27=74', '2+77=79'),
    @('35-33=2', '8+42=50'),
    @('93-11=82', '74+16=90'),
    @('6+47=53', '52-3=49'),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
